$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the newly reported evening weight / body-fat readings for 2025-CW06
# day 1 (row 2, date 45919).
$ws.Range("C2").Value = 101.7
$ws.Range("E2").Value = 27.6

# Day 2 (date 45920, row 3) had no readings logged at all - drop that blank
# row entirely, shifting the remaining days up so the tracker stays
# contiguous.
$ws.Rows("3").Delete()
